$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the formatting (bold, border,
# centered alignment) used by the other header cells (e.g. G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Add the corresponding data value in H2.
$ws.Range("H2").Value = 1
